# Add new columns I ("I0") and J ("IF") to the sheet, mirroring the
# style of the existing header row and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers I1="I0", J1="IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font + border + center/top alignment)
# from the existing header cell H1 onto the two new header cells so they
# pick up the same style used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows 2-24: column I and J values ---
# For every row, column J repeats column H's value; column I is 1 for all
# rows except row 15, where it instead takes the value 6 (matching J).
$values = @{
    2  = @(1, 3)
    3  = @(1, 5)
    4  = @(1, 3)
    5  = @(1, 3)
    6  = @(1, 5)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 7)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(1, 6)
    13 = @(1, 6)
    14 = @(1, 6)
    15 = @(6, 6)
    16 = @(1, 5)
    17 = @(1, 6)
    18 = @(1, 5)
    19 = @(1, 6)
    20 = @(1, 6)
    21 = @(1, 5)
    22 = @(1, 4)
    23 = @(1, 4)
    24 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
